$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update pais_compania (column I) values for rows 3-6
$ws.Range("I3").Value = 123
$ws.Range("I4").Value = 908
$ws.Range("I5").Value = 111
$ws.Range("I6").Value = 754

# Update correo_compania (column F) for row 6 to have a trailing space
$ws.Range("F6").Value = "9000571@hotmail.com "
